$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report header metadata (B4: period, B5: download timestamp) ---
$ws.Range("B4").Value = "2024-10-01 ~ 2024-10-31"
$ws.Range("B5").Value = "2024년 11월 06일 17시 32분 36초"

# --- The stats table gains one more day (2024-10-01), so the trailing blank spacer
#     row shifts from row 38 down to row 39. Copy its formatting down first. ---
$ws.Range("A38:F38").Copy($ws.Range("A39:F39"))

# --- Row 38 becomes a real data row; give it the same banded style as the other
#     even-numbered data rows (e.g. row 36) before writing its values. ---
$ws.Range("A36:F36").Copy($ws.Range("A38:F38"))

# --- Ensure the daily-stats table (A8:F38) stays Text-formatted so numeric-looking
#     values ("12", "0", ...) are stored as text, like the original inline-string cells ---
$ws.Range("A8:F38").NumberFormat = "@"

# --- Rewrite the daily rows for the new period (2024-10-01 .. 2024-10-31), newest first ---
$data = @(
    @("2024-10-31","목","12","0","0","12"),
    @("2024-10-30","수","11","1","0","10"),
    @("2024-10-29","화","21","0","0","21"),
    @("2024-10-28","월","12","0","0","12"),
    @("2024-10-27","일","13","0","0","13"),
    @("2024-10-26","토","9","0","0","9"),
    @("2024-10-25","금","6","0","0","6"),
    @("2024-10-24","목","4","0","0","4"),
    @("2024-10-23","수","8","0","0","8"),
    @("2024-10-22","화","13","0","0","13"),
    @("2024-10-21","월","9","0","0","9"),
    @("2024-10-20","일","12","0","0","12"),
    @("2024-10-19","토","12","0","0","12"),
    @("2024-10-18","금","12","0","0","12"),
    @("2024-10-17","목","9","0","0","9"),
    @("2024-10-16","수","18","0","0","18"),
    @("2024-10-15","화","4","1","0","3"),
    @("2024-10-14","월","16","0","0","16"),
    @("2024-10-13","일","23","3","0","20"),
    @("2024-10-12","토","30","1","1","28"),
    @("2024-10-11","금","34","3","0","31"),
    @("2024-10-10","목","19","1","0","18"),
    @("2024-10-09","수","27","6","1","20"),
    @("2024-10-08","화","18","0","0","18"),
    @("2024-10-07","월","11","0","0","11"),
    @("2024-10-06","일","14","0","0","14"),
    @("2024-10-05","토","10","0","0","10"),
    @("2024-10-04","금","15","0","0","15"),
    @("2024-10-03","목","13","1","0","12"),
    @("2024-10-02","수","21","0","0","21"),
    @("2024-10-01","화","9","0","1","8")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = 8 + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Host "done"
